$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was reported (Poroto granado, Vega Modelo de Temuco).
# It becomes the new "latest" row, right after the header, pushing the
# previously-first data row (and everything below it) down by one.
$ws.Rows(4).Insert()

$ws.Range("A4").Value = 10
$ws.Range("B4").Value = "Vega Modelo de Temuco"
$ws.Range("C4").Value = "La Araucanía"
$ws.Range("D4").Value = 44616
$ws.Range("E4").Value = 9
$ws.Range("F4").Value = 100112030
$ws.Range("G4").Value = "Poroto granado"
$ws.Range("H4").Value = "Sin especificar"
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 90
$ws.Range("K4").Value = 23000
$ws.Range("L4").Value = 25000
$ws.Range("M4").Value = 23889
$ws.Range("N4").Value = "`$/saco 25 kilos"
$ws.Range("O4").Value = "Región de La Araucanía"
$ws.Range("P4").Value = 956
$ws.Range("Q4").Value = 25
$ws.Range("R4").Value = "Hortaliza"
